$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.241.15"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "3.759.55"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.70"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").Value = "3.758.98"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  +3.57%  "
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.08"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").Value = "4.387.76"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "3.764.48"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "69.241.69"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.32"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.31"
$ws.Range("E21").Value = "  +19.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.88"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.729"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  +7.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.90"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.30"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.14"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.61"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "3.904.06"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").Value = "3.697.49"
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.02"
$ws.Range("E38").Value = "  +4.20%  "
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E40").Value = "  +2.50%  "
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("E42").Value = "  +4.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "430.94"
$ws.Range("E43").Value = "  -1.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.71"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.29"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.797.03"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.69"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("E51").Value = "  +0.44%  "
